$d = $word.ActiveDocument

# "zaposlen" + "/a" (two adjacent runs) -> single templated placeholder
# "{{ spol_zaposlen_a }}" that will later resolve to "zaposlen"/"zaposlena"
# depending on the teacher's gender (spol_zaposlen_a).
$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()
$find.Execute("zaposlen/a", $true, $false, $false, $false, $false, $true, 1, $false, "{{ spol_zaposlen_a }}", 2)
